# react 101 final scores updated
# - Project (60) scores entered for students 1, 3, 6, 9
# - "MFT" column (header + per-student formulas) removed; column F cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new Project (60) scores
$ws.Range("D15").Value = 45
$ws.Range("D17").Value = 50
$ws.Range("D20").Value = 50
$ws.Range("D23").Value = 50

# Remove the MFT column: clear the header label and the per-row formulas
$ws.Range("F14").ClearContents()
$ws.Range("F15:F25").ClearContents()

# Update the active selection to match the author's final cursor position
$ws.Range("I27").Select()
